# fix(source): fix date formatter
# Allow to user date formatter with spaces
#
# The two existing date cells (C1, C2) were stored as a bare day-count
# (43471 == 2019-01-06 00:00) and are corrected to carry the time-of-day
# fraction that the formatter now preserves (43471.6875 == 2019-01-06 16:30).
# Two brand new number formats are introduced (a datetime and a date-only
# mask using dots) plus a percentage format, each demoed on a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the existing date values so they carry the time-of-day part ---
$ws.Range("C1").Value = 43471.6875
$ws.Range("C2").Value = 43471.6875

# --- Column C needed to grow to fit the longer datetime format ---
$ws.Columns("C").ColumnWidth = 32.830729166666664

# --- New demo row: percentage formatter with a space in the format arg ---
$ws.Range("B3").Value = 3.14159
$ws.Range("B3").NumberFormat = "0%"

# --- New demo row: date formatter with a space, dd.mm.yyyy hh:mm ---
$ws.Range("C3").Value = 43471.6875
$ws.Range("C3").NumberFormat = "dd.mm.yyyy hh:mm"

# --- New demo row: date formatter with a space, dd.mm.yyyy ---
$ws.Range("C4").Value = 43471.6875
$ws.Range("C4").NumberFormat = "dd.mm.yyyy"

# --- Selection moved to reflect the newly added rows ---
[void]$ws.Range("C6").Select()
